$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1981.125
$ws.Range("J112").Value = 1980.3043
$ws.Range("L112").Value = 5940.9129
$ws.Range("N112").Value = -8156.9129
$ws.Range("H131").Value = 1504.7333
$ws.Range("I131").Value = 752.38464
$ws.Range("K131").Value = 2257.15392
$ws.Range("M131").Value = 2782.84608
$ws.Range("H137").Value = 2124.0645
$ws.Range("I137").Value = 1487.5714
$ws.Range("J137").Value = 2648.2354
$ws.Range("K137").Value = 4462.7142
$ws.Range("L137").Value = 7944.706200000001
$ws.Range("M137").Value = -1912.7142
$ws.Range("N137").Value = -13044.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 224.83333
$ws.Range("I110").Value = 204.5
$ws.Range("K110").Value = 204.5
$ws.Range("M110").Value = 1840.5
$ws.Range("H122").Value = 2573.4
$ws.Range("I122").Value = 2025.0834
$ws.Range("K122").Value = 6075.2502
$ws.Range("M122").Value = -3625.2502
$ws.Range("H132").Value = 1223.1794
$ws.Range("I132").Value = 982.82355
$ws.Range("K132").Value = 2948.47065
$ws.Range("M132").Value = -418.4706499999998
$ws.Range("H32").Value = 4972.449
$ws.Range("I32").Value = 4052.5854
$ws.Range("J32").Value = 9686.75
$ws.Range("K32").Value = 4052.5854
$ws.Range("L32").Value = 9686.75
$ws.Range("M32").Value = -3765.5854
$ws.Range("N32").Value = -10260.75
$ws.Range("H74").Value = 2865.842
$ws.Range("I74").Value = 2557
$ws.Range("K74").Value = 2557
$ws.Range("M74").Value = -1683
$ws.Range("H77").Value = 2865.842
$ws.Range("I77").Value = 2557
$ws.Range("K77").Value = 12785
$ws.Range("M77").Value = -8417

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6530.6816
$ws.Range("I134").Value = 6530.6816
$ws.Range("K134").Value = 19592.0448
$ws.Range("M134").Value = -17057.0448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1026.4615
$ws.Range("I105").Value = 1030.2727
$ws.Range("K105").Value = 1030.2727
$ws.Range("M105").Value = 716.7273
$ws.Range("H122").Value = 4781.727
$ws.Range("I122").Value = 3726.5715
$ws.Range("J122").Value = 6628.25
$ws.Range("K122").Value = 11179.7145
$ws.Range("L122").Value = 19884.75
$ws.Range("M122").Value = -8729.7145
$ws.Range("N122").Value = -24784.75
$ws.Range("H126").Value = 3534
$ws.Range("I126").Value = 1852
$ws.Range("K126").Value = 5556
$ws.Range("M126").Value = -3086
$ws.Range("H31").Value = 2395.7778
$ws.Range("I31").Value = 2060.6667
$ws.Range("K31").Value = 2060.6667
$ws.Range("M31").Value = -1765.6667
$ws.Range("H34").Value = 2395.7778
$ws.Range("I34").Value = 2060.6667
$ws.Range("K34").Value = 2060.6667
$ws.Range("M34").Value = -1858.6667
$ws.Range("H99").Value = 3534
$ws.Range("I99").Value = 1852
$ws.Range("K99").Value = 1852
$ws.Range("M99").Value = -354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 5500
$ws.Range("I123").Value = 750
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 2250
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = 200
$ws.Range("N123").Value = -49900
$ws.Range("H131").Value = 39449.777
$ws.Range("J131").Value = 50516.855
$ws.Range("L131").Value = 151550.565
$ws.Range("N131").Value = -161630.565
$ws.Range("H135").Value = 482.16666
$ws.Range("I135").Value = 398.7
$ws.Range("K135").Value = 3588.3
$ws.Range("M135").Value = -1053.3
$ws.Range("H5").Value = 482.16666
$ws.Range("I5").Value = 398.7
$ws.Range("K5").Value = 1196.1
$ws.Range("M5").Value = -1084.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4215.2085
$ws.Range("I102").Value = 5484
$ws.Range("K102").Value = 5484
$ws.Range("M102").Value = -3862
$ws.Range("H126").Value = 25746.35
$ws.Range("I126").Value = 2781.0667
$ws.Range("J126").Value = 38049.18
$ws.Range("K126").Value = 8343.2001
$ws.Range("L126").Value = 114147.54
$ws.Range("M126").Value = -5873.2001
$ws.Range("N126").Value = -119087.54
$ws.Range("H132").Value = 2380.7144
$ws.Range("I132").Value = 1693.0714
$ws.Range("J132").Value = 3756
$ws.Range("K132").Value = 5079.2142
$ws.Range("L132").Value = 11268
$ws.Range("M132").Value = -2549.2142
$ws.Range("N132").Value = -16328
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H43").Value = 1674066.5
$ws.Range("I43").Value = 2501599.8
$ws.Range("K43").Value = 2501599.8
$ws.Range("M43").Value = -2501448.8
$ws.Range("H57").Value = 20000
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1500
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3496
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 7500
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -17484
$ws.Range("H97").Value = 2121.1538
$ws.Range("I97").Value = 2186.5557
$ws.Range("J97").Value = 1974
$ws.Range("K97").Value = 2186.5557
$ws.Range("L97").Value = 1974
$ws.Range("M97").Value = -1690.5557
$ws.Range("N97").Value = -2966

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4782.0713
$ws.Range("I122").Value = 4040.7273
$ws.Range("K122").Value = 12122.1819
$ws.Range("M122").Value = -9672.1819
$ws.Range("H126").Value = 4209
$ws.Range("I126").Value = 3465
$ws.Range("K126").Value = 10395
$ws.Range("M126").Value = -7925
$ws.Range("H132").Value = 4654.909
$ws.Range("I132").Value = 4566.1333
$ws.Range("J132").Value = 4845.143
$ws.Range("K132").Value = 13698.3999
$ws.Range("L132").Value = 14535.429
$ws.Range("M132").Value = -11168.3999
$ws.Range("N132").Value = -19595.429
$ws.Range("H68").Value = 2019.091
$ws.Range("I68").Value = 1586.7142
$ws.Range("J68").Value = 2775.75
$ws.Range("K68").Value = 1586.7142
$ws.Range("L68").Value = 2775.75
$ws.Range("M68").Value = -837.7141999999999
$ws.Range("N68").Value = -4273.75
$ws.Range("H7").Value = 4209
$ws.Range("I7").Value = 3465
$ws.Range("K7").Value = 3465
$ws.Range("M7").Value = -3353
$ws.Range("H71").Value = 2019.091
$ws.Range("I71").Value = 1586.7142
$ws.Range("J71").Value = 2775.75
$ws.Range("K71").Value = 7933.571
$ws.Range("L71").Value = 13878.75
$ws.Range("M71").Value = -4189.571
$ws.Range("N71").Value = -21366.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 19793.428
$ws.Range("I126").Value = 27938.5
$ws.Range("J126").Value = 8933.333000000001
$ws.Range("K126").Value = 83815.5
$ws.Range("L126").Value = 26799.999
$ws.Range("M126").Value = -81345.5
$ws.Range("N126").Value = -31739.999
$ws.Range("H132").Value = 1989.7241
$ws.Range("I132").Value = 1155.2778
$ws.Range("K132").Value = 3465.8334
$ws.Range("M132").Value = -935.8334000000004
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -887
$ws.Range("H68").Value = 65000
$ws.Range("J68").Value = 65000
$ws.Range("L68").Value = 65000
$ws.Range("N68").Value = -66622
$ws.Range("H71").Value = 65000
$ws.Range("J71").Value = 65000
$ws.Range("L71").Value = 195000
$ws.Range("N71").Value = -203112
$ws.Range("H81").Value = 1502.8
$ws.Range("I81").Value = 628.5
$ws.Range("K81").Value = 1257
$ws.Range("M81").Value = -196
$ws.Range("H84").Value = 1502.8
$ws.Range("I84").Value = 628.5
$ws.Range("K84").Value = 6285
$ws.Range("M84").Value = -981
